# Update the "Estado de Cuenta" detail table: the previous account-statement
# periods (2308-2311) are replaced by the new periods in reverse order, and the
# "Valor Mora" amount that was attached to the oldest period now belongs to the
# newest period (the 67835 balance moved from period 2311 to period 2308's old
# row, i.e. row 16 now holds 2311/67835 and row 19 now holds 2308/100000).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periodo Mora column (E16:E19) - new period order, newest first
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2310"
$ws.Range("E18").Value = "2309"
$ws.Range("E19").Value = "2308"

# Valor Mora column (F16:F19) - the 67835 balance now belongs to period 2311 (row 16)
$ws.Range("F16").Value = 67835
$ws.Range("F17").Value = 100000
$ws.Range("F18").Value = 100000
$ws.Range("F19").Value = 100000
